# Apply updated "想去人数" (F column) counts to the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 1548
    "F4"  = 1008
    "F5"  = 11
    "F7"  = 2543
    "F9"  = 1588
    "F12" = 64
    "F13" = 496
    "F15" = 45
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
